# Reverse the order of comma-separated "Recorded By" entries in column G
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        $revParts = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $revParts)
        $cell.Value2 = $newVal
    }
}
